$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.230.15'
$ws.Cells.Item(2, 5).Value = '  -0.42%  '
$ws.Cells.Item(3, 4).Value = '1.631.27'
$ws.Cells.Item(3, 5).Value = '  -1.35%  '
$ws.Cells.Item(4, 5).Value = '  -0.09%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '216.08'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.89%  '
$ws.Cells.Item(6, 5).Value = '  +1.58%  '
$ws.Cells.Item(7, 5).Value = '  -0.05%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.256'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.36%  '
$ws.Cells.Item(9, 5).Value = '  -0.88%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '20.22'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.60%  '
$ws.Cells.Item(11, 5).Value = '  -0.14%  '
$ws.Cells.Item(12, 4).Value = '1.629.46'
$ws.Cells.Item(12, 5).Value = '  -1.30%  '
$ws.Cells.Item(13, 5).Value = '  -0.32%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.544'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.06%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '64.92'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -4.33%  '
$ws.Cells.Item(16, 4).Value = '27.210.54'
$ws.Cells.Item(16, 5).Value = '  -0.51%  '
$ws.Cells.Item(17, 4).Value = '0.0₃0734'
$ws.Cells.Item(17, 5).Value = '  -0.84%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '216.67'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.88%  '
$ws.Cells.Item(19, 5).Value = '  -0.04%  '
$ws.Cells.Item(20, 5).Value = '  +1.35%  '
$ws.Cells.Item(21, 5).Value = '  -1.35%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.45'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -3.95%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.12'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.18%  '
$ws.Cells.Item(24, 5).Value = '  +0.56%  '
$ws.Cells.Item(25, 5).Value = '  -0.14%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '7.30'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.47%  '
$ws.Cells.Item(27, 5).Value = '  -0.67%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '15.58'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -1.75%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.0506'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.59%  '
$ws.Cells.Item(30, 5).Value = '  -0.99%  '
$ws.Cells.Item(31, 5).Value = '  -0.09%  '
$ws.Cells.Item(32, 5).Value = '  -1.02%  '
$ws.Cells.Item(33, 4).Value = '1.317.80'
$ws.Cells.Item(33, 5).Value = '  +4.68%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.56'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -2.17%  '
$ws.Cells.Item(35, 5).Value = '  -0.18%  '
$ws.Cells.Item(36, 5).Value = '  -1.80%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.540'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.99%  '
$ws.Cells.Item(38, 5).Value = '  +0.03%  '
$ws.Cells.Item(39, 5).Value = '  +0.01%  '
$ws.Cells.Item(40, 5).Value = '  +0.82%  '
$ws.Cells.Item(41, 5).Value = '  -1.35%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '63.62'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +2.34%  '
$ws.Cells.Item(43, 4).Value = '1.768.75'
$ws.Cells.Item(43, 5).Value = '  -1.55%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.19'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -4.52%  '
$ws.Cells.Item(45, 5).Value = '  -1.16%  '
$ws.Cells.Item(46, 5).Value = '  -0.14%  '
$ws.Cells.Item(47, 4).Value = '0.0₆0106'
$ws.Cells.Item(47, 5).Value = '  -0.81%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.815'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +21.92%  '
$ws.Cells.Item(49, 5).Value = '  +0.24%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.55'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.76%  '
$ws.Cells.Item(51, 5).Value = '  -1.85%  '
